$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value would otherwise be auto-detected as a number
# by Excel (e.g. "604.39") are forced to Text format first so they are stored
# the same way as the other (already-text) cells in column D.
$numericPriceCells = @("D5","D6","D10","D11","D12","D13","D16","D18","D20","D21","D22","D23","D24","D25","D26","D27","D30","D31","D33","D34","D35","D37","D38","D40","D41","D42","D44","D46","D47","D48","D49","D50","D51")
foreach ($addr in $numericPriceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Column D (Price) updates
$ws.Range("D2").Value = "67.756.07"
$ws.Range("D3").Value = "3.811.84"
$ws.Range("D5").Value = "604.39"
$ws.Range("D6").Value = "167.11"
$ws.Range("D10").Value = "6.30"
$ws.Range("D11").Value = "0.451"
$ws.Range("D12").Value = "0.0000254"
$ws.Range("D13").Value = "36.04"
$ws.Range("D14").Value = "4.455.64"
$ws.Range("D15").Value = "3.814.83"
$ws.Range("D16").Value = "18.47"
$ws.Range("D17").Value = "67.805.69"
$ws.Range("D18").Value = "7.09"
$ws.Range("D20").Value = "462.91"
$ws.Range("D21").Value = "9.91"
$ws.Range("D22").Value = "0.701"
$ws.Range("D23").Value = "0.0000149"
$ws.Range("D24").Value = "83.40"
$ws.Range("D25").Value = "12.06"
$ws.Range("D26").Value = "2.11"
$ws.Range("D27").Value = "10.05"
$ws.Range("D29").Value = "3.961.73"
$ws.Range("D30").Value = "2.80"
$ws.Range("D31").Value = "7.41"
$ws.Range("D33").Value = "29.61"
$ws.Range("D34").Value = "1.00"
$ws.Range("D35").Value = "9.08"
$ws.Range("D36").Value = "3.757.45"
$ws.Range("D37").Value = "0.0999"
$ws.Range("D38").Value = "3.37"
$ws.Range("D40").Value = "0.999"
$ws.Range("D41").Value = "5.80"
$ws.Range("D42").Value = "1.00"
$ws.Range("D44").Value = "48.11"
$ws.Range("D46").Value = "28.34"
$ws.Range("D47").Value = "43.09"
$ws.Range("D48").Value = "1.40"
$ws.Range("D49").Value = "8.34"
$ws.Range("D50").Value = "148.61"
$ws.Range("D51").Value = "1.84"

# Column E (Volume 1h) updates
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("E9").Value = "  +0.43%  "
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("E22").Value = "  +0.49%  "
$ws.Range("E23").Value = "  -4.44%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("E26").Value = "  -1.34%  "
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("E29").Value = "  +0.68%  "
$ws.Range("E30").Value = "  +0.44%  "
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("E35").Value = "  -1.66%  "
$ws.Range("E36").Value = "  +0.38%  "
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("E44").Value = "  +1.99%  "
$ws.Range("E45").Value = "  +0.44%  "
$ws.Range("E46").Value = "  +10.13%  "
$ws.Range("E47").Value = "  -5.77%  "
$ws.Range("E48").Value = "  +12.00%  "
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("E51").Value = "  +0.13%  "
